$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column E (environ) with 1 for rows 2 through 9
$ws.Range("E2:E9").Value = 1

# Update the active selection to E10, matching the post-edit cursor position
$ws.Range("E10").Select()
